$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.281.26'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.904.11'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.731'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '255.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.21%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.59'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.369'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.77'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0759'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0988'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '2.179.12'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.725'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.03%  '
$ws.Range('D17').Value = '1.897.75'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '35.250.88'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.68%  '
$ws.Range('D20').Value = '0.0₃0848'
$ws.Range('E20').Value = '  +3.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.92%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.47%  '
$ws.Range('E26').Value = '  +4.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.77'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.133'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.79%  '
$ws.Range('D31').Value = '4.128.98'
$ws.Range('E31').Value = '  +19.46%  '
$ws.Range('E32').Value = '  +6.22%  '
$ws.Range('E33').Value = '  +14.74%  '
$ws.Range('E34').Value = '  +23.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0588'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.54%  '
$ws.Range('E36').Value = '  +4.53%  '
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.912'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('E40').Value = '  +5.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.10'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.21%  '
$ws.Range('E43').Value = '  +2.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0649'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.71%  '
$ws.Range('D45').Value = '1.335.34'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('E46').Value = '  +3.41%  '
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0755'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.10%  '
